$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Experience" / "Job Type" values in F2 and G2 were split apart:
# F2 used to hold the combined experience string; now it only holds the
# single experience value, and G2 holds the job type "full-time".
$ws.Range("F2").Value = "0 - 1 an experienta"
$ws.Range("G2").Value = "full-time"

# Update the active selection to reflect the author's final cursor position.
$ws.Range("C10").Select()
